# Apply the "Add files via upload" edit:
#   - rename the "Event" sheet to "Date"
#   - scroll the "Dialogue" sheet so row 96 is the top-left visible cell
#     (zoom stays at 85%)
#   - on the (still active) "Date" sheet, change the zoom to 70% and
#     reset the scroll position back to the top-left (A1)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Dialogue"
$ws2 = $wb.Worksheets.Item(2)   # "Event" -> "Date"

# Rename the second sheet.
$ws2.Name = "Date"

# Update the "Dialogue" sheet's view: scroll to A96, keep the 85% zoom.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow    = 96
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom         = 85

# Switch back to "Date" (it remains the selected/active tab), set zoom to
# 70% and scroll back to the top-left corner.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow    = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom         = 70
